$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.733.25'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '1.752.31'
$ws.Range('E3').Value = '  -4.20%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.30'
$ws.Range('E5').Value = '  -4.43%  '
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5062'
$ws.Range('E7').Value = '  -2.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.37'
$ws.Range('E8').Value = '  -6.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2636'
$ws.Range('E9').Value = '  -2.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06170'
$ws.Range('E10').Value = '  -8.15%  '
$ws.Range('D11').Value = '1.753.04'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06916'
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.60'
$ws.Range('E13').Value = '  -2.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5985'
$ws.Range('E14').Value = '  -8.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.484'
$ws.Range('E15').Value = '  -5.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '77.05'
$ws.Range('E16').Value = '  -9.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '25.750.70'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006832'
$ws.Range('E20').Value = '  -4.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.64'
$ws.Range('E21').Value = '  -9.96%  '
$ws.Range('D22').Value = '1.975.64'
$ws.Range('E22').Value = '  -4.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.079'
$ws.Range('E23').Value = '  -7.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.262'
$ws.Range('E24').Value = '  -6.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.201'
$ws.Range('E25').Value = '  -9.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '137.04'
$ws.Range('E26').Value = '  -2.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.463'
$ws.Range('E27').Value = '  -11.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.818'
$ws.Range('E28').Value = '  -8.84%  '
$ws.Range('E29').Value = '  -8.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '102.38'
$ws.Range('E30').Value = '  -4.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08197'
$ws.Range('E31').Value = '  -5.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.678'
$ws.Range('E32').Value = '  -8.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.438'
$ws.Range('E33').Value = '  -6.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04486'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.001'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.659'
$ws.Range('E36').Value = '  -7.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9910'
$ws.Range('E37').Value = '  -8.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6029'
$ws.Range('E38').Value = '  -11.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.694'
$ws.Range('E39').Value = '  -11.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01556'
$ws.Range('E40').Value = '  -3.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.946'
$ws.Range('E41').Value = '  -8.38%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '102.89'
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3800'
$ws.Range('E44').Value = '  -12.61%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7392'
$ws.Range('E45').Value = '  -12.54%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.920'
$ws.Range('E46').Value = '  -13.49%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05473'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1101'
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.932'
$ws.Range('E49').Value = '  -13.63%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.661'
$ws.Range('E50').Value = '  -9.82%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '29.88'
$ws.Range('E51').Value = '  -8.55%  '
